# "Modificacion para el dashboard"
# Update a handful of label cells on the parameters sheet and move the
# active selection, matching the authored dashboard-facing relabeling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename the two "::"-delimited combo labels and fix the
# double-dot suffixes on the office/client labels (".." -> "...").
$ws.Range("T2").Value  = "Oficina de Ventas...4"
$ws.Range("U2").Value  = "Cliente Destinatario...30"
$ws.Range("AA2").Value = "Pedido SAP::Código SAP"
$ws.Range("AB2").Value = "Pedido::Material"
$ws.Range("AD2").Value = "Pedido SAP::Código SAP"
$ws.Range("AE2").Value = "Pedido::Material"

# Row 3 / Row 4: populate the two previously-empty helper cells.
$ws.Range("M3").Value = "Ctd. Ped."
$ws.Range("M4").Value = "Cajas Naturales"

# Move the sheet's active selection to U2 (was AB5), matching the
# scrolled-into-view state saved with the workbook.
$ws.Range("U2").Select()
